# Rename the "Depth.*" stats columns to "Rel.Depth.*" and recompute the
# data in those columns as depth values relative to the max depth
# (Depth.max), per commit message "Changes (Rel.Position, Rel.Depth) to
# script and ReadMe."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# --- Header renames (O1:R1 keep Depth.n / Depth.sd as-is in N1/S1) ---
$ws.Range("O1").Value() = "Rel.Depth.min"
$ws.Range("P1").Value() = "Rel.Depth.max"
$ws.Range("Q1").Value() = "Rel.Depth.mean"
$ws.Range("R1").Value() = "Rel.Depth.median"

# --- Recompute O:R for each data row as (value - Depth.max) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $depthMin    = $ws.Cells.Item($r, 15).Value()   # column O
    $depthMax    = $ws.Cells.Item($r, 16).Value()   # column P
    $depthMean   = $ws.Cells.Item($r, 17).Value()   # column Q
    $depthMedian = $ws.Cells.Item($r, 18).Value()   # column R

    $ws.Cells.Item($r, 15).Value() = $depthMin - $depthMax
    $ws.Cells.Item($r, 16).Value() = $depthMax - $depthMax
    $ws.Cells.Item($r, 17).Value() = $depthMean - $depthMax
    $ws.Cells.Item($r, 18).Value() = $depthMedian - $depthMax
}
